$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.810.91"
$ws.Range("E2").Value = "  +8.31%  "
$ws.Range("D3").Value = "2.588.94"
$ws.Range("E3").Value = "  +10.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.00"
$ws.Range("E5").Value = "  +6.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.44"
$ws.Range("E6").Value = "  +7.91%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "2.579.21"
$ws.Range("E9").Value = "  +9.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.12"
$ws.Range("E10").Value = "  +12.42%  "
$ws.Range("E11").Value = "  +7.26%  "
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "3.018.50"
$ws.Range("E14").Value = "  +9.66%  "
$ws.Range("D15").Value = "59.584.68"
$ws.Range("E15").Value = "  +8.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.89"
$ws.Range("E16").Value = "  +9.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +6.66%  "
$ws.Range("D18").Value = "2.574.75"
$ws.Range("E18").Value = "  +9.77%  "
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.41"
$ws.Range("E20").Value = "  +8.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.39"
$ws.Range("E21").Value = "  +8.49%  "
$ws.Range("E22").Value = "  +7.90%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.12"
$ws.Range("E24").Value = "  +5.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").Value = "  +6.21%  "
$ws.Range("E26").Value = "  +8.89%  "
$ws.Range("D27").Value = "2.685.70"
$ws.Range("E27").Value = "  +9.81%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "0.0₃0834"
$ws.Range("E29").Value = "  +12.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.91"
$ws.Range("E32").Value = "  +7.88%  "
$ws.Range("E33").Value = "  +7.07%  "
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.54"
$ws.Range("E35").Value = "  +8.78%  "
$ws.Range("E36").Value = "  +9.60%  "
$ws.Range("E37").Value = "  +9.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.863"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "305.83"
$ws.Range("E39").Value = "  +21.82%  "
$ws.Range("E40").Value = "  +10.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.45"
$ws.Range("E41").Value = "  +9.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.23"
$ws.Range("E42").Value = "  +4.71%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +9.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0571"
$ws.Range("E45").Value = "  +10.59%  "
$ws.Range("E46").Value = "  +25.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").Value = "  +13.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.34"
$ws.Range("E49").Value = "  +16.14%  "
$ws.Range("E50").Value = "  +7.48%  "
$ws.Range("E51").Value = "  +0.93%  "
